# Actualización automática 2025-07-24 14:50:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D8").Value = 566.78
$ws1.Range("M8").Value = 3724.65
$ws1.Range("D24").Value = "2 de 22"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F8").Value = 4423.73
$ws2.Range("F24").Value = 30721.26

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 2032.12
$ws3.Range("E3").Value = 6636.79
$ws3.Range("F3").Value = 0.2344147072700028

$ws3.Range("D16").Value = 27184.28
$ws3.Range("E16").Value = 11572.26
$ws3.Range("F16").Value = 0.7014114263037928

$ws3.Range("D19").Value = 30721.26
$ws3.Range("E19").Value = 27501.74386304603
$ws3.Range("F19").Value = 0.5276481452633998
